$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.780.16"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.643.62"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.80"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0628"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "1.645.95"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.64"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "26.783.86"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "0.0₃0735"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.22"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.43"
$ws.Range("E21").Value = "  +7.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.30"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "146.03"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0508"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("D33").Value = "1.283.57"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.536"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.819"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.805"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.31"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").Value = "1.784.87"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.41"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.84"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.65"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0968"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +0.07%  "
